$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 193; existing rows 193-258 shift down to 194-259.
$ws.Rows("193").Insert()

# Populate the newly inserted row 193 with the new record.
$ws.Range("A193").Value = 1
$ws.Range("B193").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C193").Value = "Arica y Parinacota"
$ws.Range("D193").Value = 44755
$ws.Range("E193").Value = 15
$ws.Range("F193").Value = "Fruta"
$ws.Range("G193").Value = 100102
$ws.Range("H193").Value = "Cítricos"
$ws.Range("I193").Value = 100102003
$ws.Range("J193").Value = "Limón"
$ws.Range("K193").Value = "Sin especificar"
$ws.Range("L193").Value = "Segunda"
$ws.Range("M193").Value = 300
$ws.Range("N193").Value = 9000
$ws.Range("O193").Value = 10000
$ws.Range("P193").Value = 9500
$ws.Range("Q193").Value = "$/caja 20 kilos"
$ws.Range("R193").Value = "Región de Coquimbo"
$ws.Range("S193").Value = 475
$ws.Range("T193").Value = 20
